$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.34128212928772
$ws.Range("B1").Value = 1.600281238555908
$ws.Range("C1").Value = 4.056238174438477
$ws.Range("D1").Value = 3.20220422744751
$ws.Range("E1").Value = 1.098598480224609
